$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected/swapped rows (match rows were recorded with swapped team/market data) ---
# Row 104
$ws.Cells.Item(104, 2).Value = 7331148
$ws.Cells.Item(104, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(104, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(104, 6).Value = 'Always Ready'
$ws.Cells.Item(104, 7).Value = 'Real Santa Cruz'
$ws.Cells.Item(104, 8).Value = 2
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 'H'
$ws.Cells.Item(104, 11).Value = 1.2
$ws.Cells.Item(104, 12).Value = 6
$ws.Cells.Item(104, 13).Value = 9
$ws.Cells.Item(104, 14).Value = 1.166
$ws.Cells.Item(104, 15).Value = 7.5
$ws.Cells.Item(104, 16).Value = 16
$ws.Cells.Item(104, 17).Value = -2
$ws.Cells.Item(104, 18).Value = 1.8
$ws.Cells.Item(104, 19).Value = 2
$ws.Cells.Item(104, 20).Value = 3.5
$ws.Cells.Item(104, 21).Value = 1.975
$ws.Cells.Item(104, 22).Value = 1.825
$ws.Cells.Item(104, 23).Value = 0.1659999999999999
$ws.Cells.Item(104, 24).Value = -1
$ws.Cells.Item(104, 25).Value = -1
$ws.Cells.Item(104, 26).Value = 0
$ws.Cells.Item(104, 27).Value = -0
$ws.Cells.Item(104, 28).Value = -1
$ws.Cells.Item(104, 29).Value = 0.825

# Row 105
$ws.Cells.Item(105, 2).Value = 7331149
$ws.Cells.Item(105, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(105, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(105, 6).Value = 'Vaca Diez'
$ws.Cells.Item(105, 7).Value = 'Guabira'
$ws.Cells.Item(105, 8).Value = 3
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 'H'
$ws.Cells.Item(105, 11).Value = 1.8
$ws.Cells.Item(105, 12).Value = 3.6
$ws.Cells.Item(105, 13).Value = 3.6
$ws.Cells.Item(105, 14).Value = 2.3
$ws.Cells.Item(105, 15).Value = 3.75
$ws.Cells.Item(105, 16).Value = 2.8
$ws.Cells.Item(105, 17).Value = 0
$ws.Cells.Item(105, 18).Value = 1.8
$ws.Cells.Item(105, 19).Value = 2
$ws.Cells.Item(105, 20).Value = 2.75
$ws.Cells.Item(105, 21).Value = 1.95
$ws.Cells.Item(105, 22).Value = 1.85
$ws.Cells.Item(105, 23).Value = 1.3
$ws.Cells.Item(105, 24).Value = -1
$ws.Cells.Item(105, 25).Value = -1
$ws.Cells.Item(105, 26).Value = 0.8
$ws.Cells.Item(105, 27).Value = -1
$ws.Cells.Item(105, 28).Value = 0.475
$ws.Cells.Item(105, 29).Value = -0.5

# Row 128
$ws.Cells.Item(128, 2).Value = 7462542
$ws.Cells.Item(128, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(128, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(128, 6).Value = 'Always Ready'
$ws.Cells.Item(128, 7).Value = 'Royal Pari FC'
$ws.Cells.Item(128, 8).Value = 3
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 'H'
$ws.Cells.Item(128, 11).Value = 1.363
$ws.Cells.Item(128, 12).Value = 4.75
$ws.Cells.Item(128, 13).Value = 6.5
$ws.Cells.Item(128, 14).Value = 1.285
$ws.Cells.Item(128, 15).Value = 6.5
$ws.Cells.Item(128, 16).Value = 8
$ws.Cells.Item(128, 17).Value = -1.75
$ws.Cells.Item(128, 18).Value = 1.9
$ws.Cells.Item(128, 19).Value = 1.9
$ws.Cells.Item(128, 20).Value = 3.25
$ws.Cells.Item(128, 21).Value = 1.85
$ws.Cells.Item(128, 22).Value = 1.95
$ws.Cells.Item(128, 23).Value = 0.2849999999999999
$ws.Cells.Item(128, 24).Value = -1
$ws.Cells.Item(128, 25).Value = -1
$ws.Cells.Item(128, 26).Value = 0.8999999999999999
$ws.Cells.Item(128, 27).Value = -1
$ws.Cells.Item(128, 28).Value = -0.5
$ws.Cells.Item(128, 29).Value = 0.475

# Row 129
$ws.Cells.Item(129, 2).Value = 7462738
$ws.Cells.Item(129, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(129, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(129, 6).Value = 'Vaca Diez'
$ws.Cells.Item(129, 7).Value = 'The Strongest'
$ws.Cells.Item(129, 8).Value = 2
$ws.Cells.Item(129, 9).Value = 2
$ws.Cells.Item(129, 10).Value = 'D'
$ws.Cells.Item(129, 11).Value = 4
$ws.Cells.Item(129, 12).Value = 4
$ws.Cells.Item(129, 13).Value = 1.666
$ws.Cells.Item(129, 14).Value = 4
$ws.Cells.Item(129, 15).Value = 3.8
$ws.Cells.Item(129, 16).Value = 1.75
$ws.Cells.Item(129, 17).Value = 0.75
$ws.Cells.Item(129, 18).Value = 1.8
$ws.Cells.Item(129, 19).Value = 2
$ws.Cells.Item(129, 20).Value = 3
$ws.Cells.Item(129, 21).Value = 1.925
$ws.Cells.Item(129, 22).Value = 1.875
$ws.Cells.Item(129, 23).Value = -1
$ws.Cells.Item(129, 24).Value = 2.8
$ws.Cells.Item(129, 25).Value = -1
$ws.Cells.Item(129, 26).Value = 0.8
$ws.Cells.Item(129, 27).Value = -1
$ws.Cells.Item(129, 28).Value = 0.925
$ws.Cells.Item(129, 29).Value = -1

# Row 142
$ws.Cells.Item(142, 2).Value = 7532430
$ws.Cells.Item(142, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(142, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(142, 6).Value = 'Always Ready'
$ws.Cells.Item(142, 7).Value = 'Oriente Petrolero'
$ws.Cells.Item(142, 8).Value = 4
$ws.Cells.Item(142, 9).Value = 1
$ws.Cells.Item(142, 10).Value = 'H'
$ws.Cells.Item(142, 11).Value = 1.4
$ws.Cells.Item(142, 12).Value = 4.2
$ws.Cells.Item(142, 13).Value = 7
$ws.Cells.Item(142, 14).Value = 1.363
$ws.Cells.Item(142, 15).Value = 4.5
$ws.Cells.Item(142, 16).Value = 8.5
$ws.Cells.Item(142, 17).Value = -1.5
$ws.Cells.Item(142, 18).Value = 2
$ws.Cells.Item(142, 19).Value = 1.8
$ws.Cells.Item(142, 20).Value = 3
$ws.Cells.Item(142, 21).Value = 1.9
$ws.Cells.Item(142, 22).Value = 1.9
$ws.Cells.Item(142, 23).Value = 0.363
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = 1
$ws.Cells.Item(142, 27).Value = -1
$ws.Cells.Item(142, 28).Value = 0.8999999999999999
$ws.Cells.Item(142, 29).Value = -1

# Row 143
$ws.Cells.Item(143, 2).Value = 7532413
$ws.Cells.Item(143, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(143, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(143, 6).Value = 'Libertad Gran Mamore FC'
$ws.Cells.Item(143, 7).Value = 'Club Aurora'
$ws.Cells.Item(143, 8).Value = 0
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 10).Value = 'A'
$ws.Cells.Item(143, 11).Value = 2.25
$ws.Cells.Item(143, 12).Value = 3.3
$ws.Cells.Item(143, 13).Value = 2.8
$ws.Cells.Item(143, 14).Value = 2.375
$ws.Cells.Item(143, 15).Value = 3.4
$ws.Cells.Item(143, 16).Value = 2.875
$ws.Cells.Item(143, 17).Value = -0.25
$ws.Cells.Item(143, 18).Value = 2.025
$ws.Cells.Item(143, 19).Value = 1.775
$ws.Cells.Item(143, 20).Value = 2.5
$ws.Cells.Item(143, 21).Value = 1.9
$ws.Cells.Item(143, 22).Value = 1.9
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 24).Value = -1
$ws.Cells.Item(143, 25).Value = 1.875
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = 0.7749999999999999
$ws.Cells.Item(143, 28).Value = -1
$ws.Cells.Item(143, 29).Value = 0.8999999999999999

# Row 144
$ws.Cells.Item(144, 2).Value = 7532414
$ws.Cells.Item(144, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(144, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(144, 6).Value = 'Independiente Petrolero'
$ws.Cells.Item(144, 7).Value = 'Real Santa Cruz'
$ws.Cells.Item(144, 8).Value = 1
$ws.Cells.Item(144, 9).Value = 0
$ws.Cells.Item(144, 10).Value = 'H'
$ws.Cells.Item(144, 11).Value = 1.571
$ws.Cells.Item(144, 12).Value = 3.75
$ws.Cells.Item(144, 13).Value = 5
$ws.Cells.Item(144, 14).Value = 1.3
$ws.Cells.Item(144, 15).Value = 5
$ws.Cells.Item(144, 16).Value = 11
$ws.Cells.Item(144, 17).Value = -1.75
$ws.Cells.Item(144, 18).Value = 2
$ws.Cells.Item(144, 19).Value = 1.8
$ws.Cells.Item(144, 20).Value = 3
$ws.Cells.Item(144, 21).Value = 1.85
$ws.Cells.Item(144, 22).Value = 1.95
$ws.Cells.Item(144, 23).Value = 0.3
$ws.Cells.Item(144, 24).Value = -1
$ws.Cells.Item(144, 25).Value = -1
$ws.Cells.Item(144, 26).Value = -1
$ws.Cells.Item(144, 27).Value = 0.8
$ws.Cells.Item(144, 28).Value = -1
$ws.Cells.Item(144, 29).Value = 0.95

# Row 148
$ws.Cells.Item(148, 2).Value = 7532421
$ws.Cells.Item(148, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(148, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(148, 6).Value = 'Guabira'
$ws.Cells.Item(148, 7).Value = 'Independiente Petrolero'
$ws.Cells.Item(148, 8).Value = 2
$ws.Cells.Item(148, 9).Value = 0
$ws.Cells.Item(148, 10).Value = 'H'
$ws.Cells.Item(148, 11).Value = 1.4
$ws.Cells.Item(148, 12).Value = 4.5
$ws.Cells.Item(148, 13).Value = 7.5
$ws.Cells.Item(148, 14).Value = 1.333
$ws.Cells.Item(148, 15).Value = 5.5
$ws.Cells.Item(148, 16).Value = 9.5
$ws.Cells.Item(148, 17).Value = -1.5
$ws.Cells.Item(148, 18).Value = 1.85
$ws.Cells.Item(148, 19).Value = 1.95
$ws.Cells.Item(148, 20).Value = 3
$ws.Cells.Item(148, 21).Value = 1.825
$ws.Cells.Item(148, 22).Value = 1.975
$ws.Cells.Item(148, 23).Value = 0.333
$ws.Cells.Item(148, 24).Value = -1
$ws.Cells.Item(148, 25).Value = -1
$ws.Cells.Item(148, 26).Value = 0.8500000000000001
$ws.Cells.Item(148, 27).Value = -1
$ws.Cells.Item(148, 28).Value = -1
$ws.Cells.Item(148, 29).Value = 0.9750000000000001

# Row 150
$ws.Cells.Item(150, 2).Value = 7532419
$ws.Cells.Item(150, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(150, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(150, 6).Value = 'Oriente Petrolero'
$ws.Cells.Item(150, 7).Value = 'Jorge Wilstermann'
$ws.Cells.Item(150, 8).Value = 3
$ws.Cells.Item(150, 9).Value = 0
$ws.Cells.Item(150, 10).Value = 'H'
$ws.Cells.Item(150, 11).Value = 2.2
$ws.Cells.Item(150, 12).Value = 2.5
$ws.Cells.Item(150, 13).Value = 4.5
$ws.Cells.Item(150, 14).Value = 2.375
$ws.Cells.Item(150, 15).Value = 2.45
$ws.Cells.Item(150, 16).Value = 4.5
$ws.Cells.Item(150, 17).Value = -0.25
$ws.Cells.Item(150, 18).Value = 1.9
$ws.Cells.Item(150, 19).Value = 1.9
$ws.Cells.Item(150, 20).Value = 2
$ws.Cells.Item(150, 21).Value = 1.95
$ws.Cells.Item(150, 22).Value = 1.85
$ws.Cells.Item(150, 23).Value = 1.375
$ws.Cells.Item(150, 24).Value = -1
$ws.Cells.Item(150, 25).Value = -1
$ws.Cells.Item(150, 26).Value = 0.8999999999999999
$ws.Cells.Item(150, 27).Value = -1
$ws.Cells.Item(150, 28).Value = 0.95
$ws.Cells.Item(150, 29).Value = -1

# Row 153
$ws.Cells.Item(153, 2).Value = 7532417
$ws.Cells.Item(153, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(153, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(153, 6).Value = 'Real Tomayapo'
$ws.Cells.Item(153, 7).Value = 'Atletico Palmaflor Vinto'
$ws.Cells.Item(153, 8).Value = 4
$ws.Cells.Item(153, 9).Value = 0
$ws.Cells.Item(153, 10).Value = 'H'
$ws.Cells.Item(153, 11).Value = 1.3
$ws.Cells.Item(153, 12).Value = 4.5
$ws.Cells.Item(153, 13).Value = 8.5
$ws.Cells.Item(153, 14).Value = 1.166
$ws.Cells.Item(153, 15).Value = 8
$ws.Cells.Item(153, 16).Value = 12
$ws.Cells.Item(153, 17).Value = -2.25
$ws.Cells.Item(153, 18).Value = 1.95
$ws.Cells.Item(153, 19).Value = 1.85
$ws.Cells.Item(153, 20).Value = 3.75
$ws.Cells.Item(153, 21).Value = 1.975
$ws.Cells.Item(153, 22).Value = 1.825
$ws.Cells.Item(153, 23).Value = 0.1659999999999999
$ws.Cells.Item(153, 24).Value = -1
$ws.Cells.Item(153, 25).Value = -1
$ws.Cells.Item(153, 26).Value = 0.95
$ws.Cells.Item(153, 27).Value = -1
$ws.Cells.Item(153, 28).Value = 0.4875
$ws.Cells.Item(153, 29).Value = -0.5

# Row 154
$ws.Cells.Item(154, 2).Value = 7532431
$ws.Cells.Item(154, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(154, 4).Value = 'Bolivia Clausura'
$ws.Cells.Item(154, 6).Value = 'Blooming'
$ws.Cells.Item(154, 7).Value = 'Always Ready'
$ws.Cells.Item(154, 8).Value = 2
$ws.Cells.Item(154, 9).Value = 1
$ws.Cells.Item(154, 10).Value = 'H'
$ws.Cells.Item(154, 11).Value = 2.3
$ws.Cells.Item(154, 12).Value = 3.5
$ws.Cells.Item(154, 13).Value = 2.625
$ws.Cells.Item(154, 14).Value = 1.833
$ws.Cells.Item(154, 15).Value = 4
$ws.Cells.Item(154, 16).Value = 3.6
$ws.Cells.Item(154, 17).Value = -0.5
$ws.Cells.Item(154, 18).Value = 1.825
$ws.Cells.Item(154, 19).Value = 1.975
$ws.Cells.Item(154, 20).Value = 3
$ws.Cells.Item(154, 21).Value = 2
$ws.Cells.Item(154, 22).Value = 1.8
$ws.Cells.Item(154, 23).Value = 0.833
$ws.Cells.Item(154, 24).Value = -1
$ws.Cells.Item(154, 25).Value = -1
$ws.Cells.Item(154, 26).Value = 0.825
$ws.Cells.Item(154, 27).Value = -1
$ws.Cells.Item(154, 28).Value = 0
$ws.Cells.Item(154, 29).Value = -0

# --- New rows appended at the end (203-206) ---
# Row 203
$ws.Cells.Item(203, 1).Value = 201
$ws.Cells.Item(203, 2).Value = 8010639
$ws.Cells.Item(203, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(203, 4).Value = 'Bolivia Apertura'
$ws.Cells.Item(203, 5).Value = 45381.875
$ws.Cells.Item(203, 6).Value = 'Nacional Potosi'
$ws.Cells.Item(203, 7).Value = 'Royal Pari FC'
$ws.Cells.Item(203, 11).Value = 1.363
$ws.Cells.Item(203, 12).Value = 4.333
$ws.Cells.Item(203, 13).Value = 7
$ws.Cells.Item(203, 14).Value = 1.333
$ws.Cells.Item(203, 15).Value = 5
$ws.Cells.Item(203, 16).Value = 9.5
$ws.Cells.Item(203, 17).Value = -1.5
$ws.Cells.Item(203, 18).Value = 1.825
$ws.Cells.Item(203, 19).Value = 1.975
$ws.Cells.Item(203, 20).Value = 3.25
$ws.Cells.Item(203, 21).Value = 1.8
$ws.Cells.Item(203, 22).Value = 2
$ws.Cells.Item(203, 23).Value = 0
$ws.Cells.Item(203, 24).Value = 0
$ws.Cells.Item(203, 25).Value = 0
$ws.Cells.Item(203, 26).Value = 0
$ws.Cells.Item(203, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Cells.Item(203, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(203, 5).PasteSpecial(-4122)

# Row 204
$ws.Cells.Item(204, 1).Value = 202
$ws.Cells.Item(204, 2).Value = 8010642
$ws.Cells.Item(204, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(204, 4).Value = 'Bolivia Apertura'
$ws.Cells.Item(204, 5).Value = 45382.66666666666
$ws.Cells.Item(204, 6).Value = 'Always Ready'
$ws.Cells.Item(204, 7).Value = 'Independiente Petrolero'
$ws.Cells.Item(204, 11).Value = 1.4
$ws.Cells.Item(204, 12).Value = 4
$ws.Cells.Item(204, 13).Value = 7
$ws.Cells.Item(204, 14).Value = 1.444
$ws.Cells.Item(204, 15).Value = 4.333
$ws.Cells.Item(204, 16).Value = 8
$ws.Cells.Item(204, 17).Value = -1.25
$ws.Cells.Item(204, 18).Value = 1.9
$ws.Cells.Item(204, 19).Value = 1.9
$ws.Cells.Item(204, 20).Value = 2.75
$ws.Cells.Item(204, 21).Value = 1.8
$ws.Cells.Item(204, 22).Value = 2
$ws.Cells.Item(204, 23).Value = 0
$ws.Cells.Item(204, 24).Value = 0
$ws.Cells.Item(204, 25).Value = 0
$ws.Cells.Item(204, 26).Value = 0
$ws.Cells.Item(204, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Cells.Item(204, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(204, 5).PasteSpecial(-4122)

# Row 205
$ws.Cells.Item(205, 1).Value = 203
$ws.Cells.Item(205, 2).Value = 8010640
$ws.Cells.Item(205, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(205, 4).Value = 'Bolivia Apertura'
$ws.Cells.Item(205, 5).Value = 45382.77083333334
$ws.Cells.Item(205, 6).Value = 'Bolivar'
$ws.Cells.Item(205, 7).Value = 'Oriente Petrolero'
$ws.Cells.Item(205, 11).Value = 1.222
$ws.Cells.Item(205, 12).Value = 5.75
$ws.Cells.Item(205, 13).Value = 9
$ws.Cells.Item(205, 14).Value = 1.222
$ws.Cells.Item(205, 15).Value = 6
$ws.Cells.Item(205, 16).Value = 13
$ws.Cells.Item(205, 17).Value = -1.75
$ws.Cells.Item(205, 18).Value = 1.8
$ws.Cells.Item(205, 19).Value = 2
$ws.Cells.Item(205, 20).Value = 3.25
$ws.Cells.Item(205, 21).Value = 2.025
$ws.Cells.Item(205, 22).Value = 1.775
$ws.Cells.Item(205, 23).Value = 0
$ws.Cells.Item(205, 24).Value = 0
$ws.Cells.Item(205, 25).Value = 0
$ws.Cells.Item(205, 26).Value = 0
$ws.Cells.Item(205, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Cells.Item(205, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(205, 5).PasteSpecial(-4122)

# Row 206
$ws.Cells.Item(206, 1).Value = 204
$ws.Cells.Item(206, 2).Value = 8011587
$ws.Cells.Item(206, 3).Value = 'Bolivia Primera División'
$ws.Cells.Item(206, 4).Value = 'Bolivia Apertura'
$ws.Cells.Item(206, 5).Value = 45382.85416666666
$ws.Cells.Item(206, 6).Value = 'Club Aurora'
$ws.Cells.Item(206, 7).Value = 'Blooming'
$ws.Cells.Item(206, 11).Value = 1.533
$ws.Cells.Item(206, 12).Value = 3.75
$ws.Cells.Item(206, 13).Value = 5.5
$ws.Cells.Item(206, 14).Value = 1.4
$ws.Cells.Item(206, 15).Value = 4.333
$ws.Cells.Item(206, 16).Value = 8.5
$ws.Cells.Item(206, 17).Value = -1.25
$ws.Cells.Item(206, 18).Value = 1.925
$ws.Cells.Item(206, 19).Value = 1.875
$ws.Cells.Item(206, 20).Value = 2.75
$ws.Cells.Item(206, 21).Value = 1.85
$ws.Cells.Item(206, 22).Value = 1.95
$ws.Cells.Item(206, 23).Value = 0
$ws.Cells.Item(206, 24).Value = 0
$ws.Cells.Item(206, 25).Value = 0
$ws.Cells.Item(206, 26).Value = 0
$ws.Cells.Item(206, 27).Value = 0
$ws.Range("A2").Copy()
$ws.Cells.Item(206, 1).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Cells.Item(206, 5).PasteSpecial(-4122)

